# Update "想去人数" (interest count) figures in column F across sheets,
# matching the regenerated gh-pages data snapshot at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 4867  # was 4860
$ws.Range("F8").Value = 4867  # was 4860
$ws.Range("F13").Value = 1111  # was 1110
$ws.Range("F15").Value = 4523  # was 4512
$ws.Range("F17").Value = 187  # was 186
$ws.Range("F18").Value = 81  # was 79
$ws.Range("F19").Value = 228  # was 227
$ws.Range("F20").Value = 3576  # was 3571
$ws.Range("F21").Value = 6  # was 5
$ws.Range("F24").Value = 3307  # was 3298
$ws.Range("F25").Value = 149  # was 148
$ws.Range("F26").Value = 138  # was 137
$ws.Range("F28").Value = 344  # was 343
$ws.Range("F30").Value = 210  # was 208
$ws.Range("F31").Value = 186  # was 184
$ws.Range("F37").Value = 5769  # was 5745
$ws.Range("F38").Value = 907  # was 903
$ws.Range("F39").Value = 427  # was 425
$ws.Range("F43").Value = 1164  # was 1159
$ws.Range("F44").Value = 532  # was 529
$ws.Range("F45").Value = 19  # was 18
$ws.Range("F46").Value = 2049  # was 2046
$ws.Range("F49").Value = 727  # was 725

# Sheet 2: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 125  # was 122
$ws.Range("F9").Value = 42  # was 41
$ws.Range("F24").Value = 760  # was 757

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 206  # was 205

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 206  # was 205
$ws.Range("F8").Value = 4867  # was 4860
$ws.Range("F9").Value = 4867  # was 4860
$ws.Range("F11").Value = 42  # was 41
$ws.Range("F15").Value = 1111  # was 1110
$ws.Range("F17").Value = 4523  # was 4513
$ws.Range("F19").Value = 187  # was 186
$ws.Range("F20").Value = 81  # was 79
$ws.Range("F21").Value = 228  # was 227
$ws.Range("F22").Value = 3576  # was 3571
$ws.Range("F23").Value = 3307  # was 3298
$ws.Range("F24").Value = 149  # was 148
$ws.Range("F25").Value = 138  # was 137
$ws.Range("F27").Value = 210  # was 208
$ws.Range("F28").Value = 186  # was 184
$ws.Range("F35").Value = 5769  # was 5745
$ws.Range("F37").Value = 907  # was 903
$ws.Range("F38").Value = 427  # was 425
$ws.Range("F44").Value = 1164  # was 1159
$ws.Range("F45").Value = 532  # was 529
$ws.Range("F46").Value = 2049  # was 2046
$ws.Range("F48").Value = 727  # was 725

